# Daily attendance processing - normalize "Recorded By" (column G) values so
# that the "System" entry is listed first instead of last, for every row
# whose recorded-by list includes a "System"/"system" token.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Row + $used.Rows.Count - 1
$colG = 7

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $colG)
    $raw = $cell.Value2

    if ($raw -eq $null) {
        continue
    }

    $text = [string]$raw
    if ($text -eq "") {
        continue
    }

    $rawParts = $text.Split(",")
    if ($rawParts.Count -le 1) {
        continue
    }

    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $hasSystem = $false
    foreach ($p in $parts) {
        if ($p.ToLower() -eq "system") {
            $hasSystem = $true
        }
    }

    if (-not $hasSystem) {
        continue
    }

    # Swap the first and last comma-separated entries so the "System" token
    # moves to the front of the list (matching the new recorded-by ordering).
    $newParts = @()
    $newParts += $parts[$parts.Count - 1]
    for ($i = 1; $i -lt ($parts.Count - 1); $i++) {
        $newParts += $parts[$i]
    }
    $newParts += $parts[0]

    $result = [string]::Join(", ", $newParts)

    # Always write back: PowerShell string comparisons here are
    # case-insensitive, which would otherwise mask swaps that only change
    # casing (e.g. "system, X, System" -> "System, X, system").
    $cell.Value2 = $result
}
